# ------------------------------------------------------------------
# "creation reu avec nouvelle base"
#
# 1) The "datetimeFigureOut" date placeholder text on the slide master
#    and every slide layout changes from 09/03/2020 to 19/03/2020.
# 2) Three shapes on slide 1 get resized/repositioned:
#      - "Ellipse 36"        : width/height grow
#      - "Connecteur droit 12": shifts down slightly, shrinks height
#      - "Connecteur droit 38": shifts right, narrower but taller
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# EMU -> point conversion constant (914400 EMU per inch, 72 points per inch)
$EMU_PER_PT = 12700
# Tiny nudge (well under half an EMU) so float round-trip through the
# point-based COM properties lands back on the exact target EMU value
# instead of being truncated one EMU short.
$EPS = 0.3 / $EMU_PER_PT

function EmuToPt([double]$emu) {
    return ($emu / $EMU_PER_PT) + $EPS
}

# --- 1) Update the date placeholder text (master + all layouts) ----
$oldDate = "09/03/2020"
$newDate = "19/03/2020"

function UpdateDatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

UpdateDatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    UpdateDatePlaceholder $layouts.Item($li).Shapes
}

# --- 2) Resize / reposition the three shapes on slide 1 ------------
$slide1 = $p.Slides.Item(1)
$shapes1 = $slide1.Shapes

for ($i = 1; $i -le $shapes1.Count; $i++) {
    $shp = $shapes1.Item($i)

    if ($shp.Name -eq "Ellipse 36") {
        $shp.Width = EmuToPt 1606828
        $shp.Height = EmuToPt 1274849
    }
    elseif ($shp.Name -eq "Connecteur droit 12") {
        $shp.Top = EmuToPt 2650034
        $shp.Height = EmuToPt 329399
    }
    elseif ($shp.Name -eq "Connecteur droit 38") {
        $shp.Left = EmuToPt 7586197
        $shp.Width = EmuToPt 609558
        $shp.Height = EmuToPt 1443814
    }
}
